$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.703.76"
$ws.Range("E2").Value = "  -5.03%  "
$ws.Range("D3").Value = "1.809.19"
$ws.Range("E3").Value = "  -4.02%  "
$ws.Range("E4").Value = "  +0.12%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "275.98"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -9.12%  "
$ws.Range("E6").Value = "  +0.13%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.5038"
$c.Style = "Normal"
$ws.Range("E7").Value = "  -5.80%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.3492"
$c.Style = "Normal"
$ws.Range("E8").Value = "  -7.88%  "
$ws.Range("B9").Value = "OKB"
$ws.Range("C9").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "44.45"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -2.84%  "
$ws.Range("B10").Value = "Dogecoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.06645"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -7.96%  "
$ws.Range("B11").Value = "Solana"
$ws.Range("C11").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "19.81"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -9.56%  "
$ws.Range("B12").Value = "Polygon"
$ws.Range("C12").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.8269"
$c.Style = "Normal"
$ws.Range("E12").Value = "  -7.19%  "
$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "0.07841"
$c.Style = "Normal"
$ws.Range("E13").Value = "  -3.68%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.802.44"
$ws.Range("E14").Value = "  -5.62%  "
$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "5.043"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -5.09%  "
$ws.Range("B16").Value = "Litecoin"
$ws.Range("C16").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "87.23"
$c.Style = "Normal"
$ws.Range("E16").Value = "  -7.02%  "
$ws.Range("B17").Value = "BinanceUSD"
$ws.Range("C17").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "1.000"
$c.Style = "Normal"
$ws.Range("E17").Value = "  +0.09%  "
$ws.Range("B18").Value = "Avalanche"
$ws.Range("C18").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "13.92"
$c.Style = "Normal"
$ws.Range("E18").Value = "  -6.23%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "0.000008030"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -7.17%  "
$ws.Range("B20").Value = "Dai"
$ws.Range("C20").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "1.001"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +0.18%  "
$ws.Range("B21").Value = "WrappedBTC"
$ws.Range("C21").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D21").Value = "25.765.50"
$ws.Range("E21").Value = "  -4.85%  "
$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "4.703"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -6.12%  "
$ws.Range("B23").Value = "Cosmos"
$ws.Range("C23").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "9.955"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -7.86%  "
$ws.Range("B24").Value = "Chainlink"
$ws.Range("C24").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "6.025"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -6.41%  "
$ws.Range("B25").Value = "LidoDAOToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "2.196"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -3.94%  "
$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "139.47"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -5.77%  "
$ws.Range("B27").Value = "Toncoin"
$ws.Range("C27").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "1.665"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -4.44%  "
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "17.00"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -6.46%  "
$ws.Range("B29").Value = "BitcoinCash"
$ws.Range("C29").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "109.24"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -5.22%  "
$ws.Range("B30").Value = "InternetComputer(DFINITY)"
$ws.Range("C30").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "4.311"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -9.53%  "
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "4.207"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -9.39%  "
$ws.Range("B32").Value = "Stellar"
$ws.Range("C32").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "0.08785"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -4.12%  "
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "0.04854"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -3.59%  "
$ws.Range("B34").Value = "ARBITRUM"
$ws.Range("C34").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "1.132"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -4.40%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "2.884"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -4.23%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "0.7131"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -11.93%  "
$ws.Range("B37").Value = "Frax"
$ws.Range("C37").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "1.000"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +0.22%  "
$ws.Range("B38").Value = "MXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "3.098"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -4.48%  "
$ws.Range("B39").Value = "TheSandbox"
$ws.Range("C39").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.5169"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -13.20%  "
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.01828"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -7.14%  "
$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "2.230"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -15.77%  "
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.9520"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -11.24%  "
$ws.Range("B43").Value = "Quant"
$ws.Range("C43").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "113.40"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -2.51%  "
$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "6.115"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -7.80%  "
$ws.Range("B45").Value = "Aptos"
$ws.Range("C45").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "8.004"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -10.82%  "
$ws.Range("B46").Value = "PaxDollar"
$ws.Range("C46").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "1.000"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +0.18%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "0.1367"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -9.53%  "
$ws.Range("B48").Value = "Decentraland"
$ws.Range("C48").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "0.4512"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -11.62%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "9.254"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -8.26%  "
$ws.Range("B50").Value = "Elrond"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "36.19"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -3.84%  "
$ws.Range("B51").Value = "NEARProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "1.493"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -7.77%  "
